# Add the Łukasz Napora "Komentarz" (column E) notes for weeks 1-5 of the
# work plan, and move the active selection to E10 (next empty comment cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "Pytania odnośnie dokumentacji i podział prac w grupie. Projektowanie biblioteki dla całego zespołu (3h)"
$ws.Range("E6").Value = "Stworzenie szkieletu serwera i dokończenie biblioteki (3h)"
$ws.Range("E7").Value = "Implementacja odbierania i wysyłania wiadomości (3h)"
$ws.Range("E8").Value = "Implementacja prostego algorytmu rozdzielającego zadania pomiędzy komponenty (3h)"
$ws.Range("E9").Value = "Wprowadzanie poprawek zauważonych podczas testów (3h)"

$ws.Range("E10").Select() | Out-Null
